$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6913084387779236
$ws.Range("B1").Value = 2.211065769195557
$ws.Range("D1").Value = 0.963458776473999
$ws.Range("E1").Value = 1.048778891563416
